$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "c4"
$ws.Range("B1").Value = 4
$ws.Range("A2").Value = "d4"
$ws.Range("B2").Value = 2
$ws.Range("A3").Value = "e4"
$ws.Range("B3").Value = 4
$ws.Range("A4").Value = "f4"
$ws.Range("B4").Value = 4
$ws.Range("A5").Value = "e4"
$ws.Range("B5").Value = 4
$ws.Range("A6").Value = "c4"
$ws.Range("B6").Value = 4
$ws.Range("A7").Value = "g3"
$ws.Range("B7").Value = 4
$ws.Range("A8").Value = "b3"
$ws.Range("B8").Value = 4
$ws.Range("A9").Value = "c4"
$ws.Range("B9").Value = 4
$ws.Range("A10").Value = "c5"
$ws.Range("B10").Value = 4
$ws.Range("A11").Value = "g4"
$ws.Range("B11").Value = 4
$ws.Range("A12").Value = "e4"
$ws.Range("B12").Value = 4
$ws.Range("A13").Value = "a4"
$ws.Range("B13").Value = 2
$ws.Range("A14").Value = "g4"
$ws.Range("B14").Value = 4
$ws.Range("A15").Value = "c5"
$ws.Range("B15").Value = 4
$ws.Range("A16").Value = "g4"
$ws.Range("B16").Value = 2
$ws.Range("A17").Value = "a4"
$ws.Range("B17").Value = 4
$ws.Range("A18").Value = "c5"
$ws.Range("B18").Value = 4
$ws.Range("A19").Value = "c5"
$ws.Range("B19").Value = 4
$ws.Range("A20").Value = "d5"
$ws.Range("B20").Value = 2
$ws.Range("A21").Value = "e5"
$ws.Range("B21").Value = 4
$ws.Range("A22").Value = "d5"
$ws.Range("B22").Value = 4
$ws.Range("A23").Value = "c5"
$ws.Range("B23").Value = 2
$ws.Range("A24").Value = "a4"
$ws.Range("B24").Value = 4
$ws.Range("A25").Value = "d5"
$ws.Range("B25").Value = 4
$ws.Range("A26").Value = "c5"
$ws.Range("B26").Value = 4
$ws.Range("A27").Value = "a4"
$ws.Range("B27").Value = 4
$ws.Range("A28").Value = "d5"
$ws.Range("B28").Value = 2
$ws.Range("A29").Value = "c5"
$ws.Range("B29").Value = 4
$ws.Range("A30").Value = "g4"
$ws.Range("B30").Value = 4
$ws.Range("A31").Value = "e4"
$ws.Range("B31").Value = 4
$ws.Range("A32").Value = "a4"
$ws.Range("B32").Value = 2
$ws.Range("A33").Value = "g4"
$ws.Range("B33").Value = 4
$ws.Range("A34").Value = "c5"
$ws.Range("B34").Value = 4
$ws.Range("A35").Value = "g4"
$ws.Range("B35").Value = 2
$ws.Range("A36").Value = "a4"
$ws.Range("B36").Value = 4
$ws.Range("A37").Value = "d4"
$ws.Range("B37").Value = 4
$ws.Range("A38").Value = "d4"
$ws.Range("B38").Value = 2
$ws.Range("A39").Value = "e4"
$ws.Range("B39").Value = 4
$ws.Range("A40").Value = "f4"
$ws.Range("B40").Value = 4
$ws.Range("A41").Value = "g4"
$ws.Range("B41").Value = 4
$ws.Range("A42").Value = "c4"
$ws.Range("B42").Value = 2
$ws.Range("A43").Value = "c4"
$ws.Range("B43").Value = 4
$ws.Range("A44").Value = "g4"
$ws.Range("B44").Value = 4
$ws.Range("A45").Value = "g4"
$ws.Range("B45").Value = 4
$ws.Range("A46").Value = "a4"
$ws.Range("B46").Value = 4
$ws.Range("A47").Value = "d4"
$ws.Range("B47").Value = 4
$ws.Range("A48").Value = "d4"
$ws.Range("B48").Value = 4
$ws.Range("A49").Value = "e4"
$ws.Range("B49").Value = 4
$ws.Range("A50").Value = "c4"
$ws.Range("B50").Value = 4
$ws.Range("A51").Value = "g3"
$ws.Range("B51").Value = 2
$ws.Range("A52").Value = "c4"
$ws.Range("B52").Value = 4
$ws.Range("A53").Value = "d4"
$ws.Range("B53").Value = 4
$ws.Range("A54").Value = "e4"
$ws.Range("B54").Value = 4
$ws.Range("A55").Value = "f4"
$ws.Range("B55").Value = 2
$ws.Range("A56").Value = "g4"
$ws.Range("B56").Value = 4
$ws.Range("A57").Value = "c4"
$ws.Range("B57").Value = 4
$ws.Range("A58").Value = "c4"
$ws.Range("B58").Value = 2
$ws.Range("A59").Value = "g3"
$ws.Range("B59").Value = 4
$ws.Range("A60").Value = "c4"
$ws.Range("B60").Value = 4
$ws.Range("A61").Value = "d4"
$ws.Range("B61").Value = 2
$ws.Range("A62").Value = "e4"
$ws.Range("B62").Value = 4
